# ProQ2-group1-BOM.xlsx — "Some more todo items and BOM components"
#
# Fill in quantities (column C) for several BOM line items that previously
# had no Amt value, bump the height of a few rows to match their
# neighbours, and move the sheet's scroll/selection down to where the
# user was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New quantities for the R1..R6 resistor rows (preamp section) -------
$ws.Range("C18").Value = 1   # R1: 100
$ws.Range("C19").Value = 1   # R2: 1M
$ws.Range("C20").Value = 1   # R3: 56k
$ws.Range("C21").Value = 1   # R5: 1k
$ws.Range("C22").Value = 1   # R6: 22k

# --- New quantities for the capacitor rows, which also grow a bit taller
$ws.Range("C24").Value = 2   # C3/C4: 100nF
$ws.Range("C25").Value = 1   # C5: 2.2uF
$ws.Range("C26").Value = 1   # C6: 22uF

$ws.Rows("24:26").RowHeight = 14.9

# --- Move the view to where editing left off -----------------------------
$excel.Goto($ws.Range("F32"))
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
